$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date and Description property values ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Update the Date property value
$wsMeta.Range("B8").Value = "2026-01-16T13:49:34+00:00"

# Update the Description property value
$wsMeta.Range("B12").Value = "Instructions au patient"

# --- Elements sheet: the root element's Definition mirrors the resource Description ---
$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("M2").Value = "Instructions au patient"
